$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.01113415480545901
$ws.Range("B3").Value = 0.8108811285997877
$ws.Range("B4").Value = 93.62568961948261
$ws.Range("B5").Value = "[1.0, 0.16587883390732558, 0.1737308605979352, 0.14316903807132492, 0.108581244783241, 0.34366576030772633, 0.12020598228464434, 0.1347600818044229, 0.10590525330660626, 0.08370637385165386, 0.06641034785321234, 0.06608861524405524, 0.08389848744021032, 0.09425934573975436, 0.1304012214362287, 0.11722249185014903, 0.08191929820863048, 0.05993516418432766, 0.087463629072358, 0.10127371886214095]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -3.284423136558633
$ws.Range("B5").Value = 0.7367233751372867
$ws.Range("B6").Value = 66.81380953773717
$ws.Range("B7").Value = 3.922480780724642
$ws.Range("B8").Value = 0.6701410476351632
$ws.Range("B9").Value = 57.95287080899256
$ws.Range("B10").Value = "[1.0, 0.2027334670527049, 0.17050869852100423, 0.14434168333543046, 0.15759880736524873, 0.23469516897431497, 0.09835014564604287, 0.09515598729730786, 0.10205561917727075, 0.060549962940321236, 0.05356550747231093, 0.06129136452416732, 0.05663819061716943, 0.06657332797219698, 0.09231776360412296, 0.06946737907457236, 0.06639410066619329, 0.0444703007967054, 0.07878369704911602, 0.07225219615215117]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.009597108611915059
$ws.Range("B3").Value = 1.741404077594689
$ws.Range("B4").Value = 0.07306292907826636
$ws.Range("B5").Value = 0.9118036367423069
$ws.Range("B6").Value = "[0.9999999999999999, -0.015403821002664768, -0.021940679004667188, -0.04913281505074665, -0.07028368780820968, 0.2903129763184927, -0.002213471025706947, 0.01610503420020718, -0.017693341677344213, -0.023382053783121907, -0.05800503861475942, -0.06130107745523617, -0.023816799036103245, -0.0005576054742043746, 0.05078913064788003, 0.010747409324063472, -0.027834730011086217, -0.057042973221021434, -0.030336916478557328, -0.004010179747157036]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.01061789751592729
$ws.Range("B3").Value = 1.739570152053577
$ws.Range("B4").Value = 0.06906479509566742
$ws.Range("B5").Value = 0.9120602775760585
$ws.Range("B6").Value = "[1.0, -0.015187448945871967, -0.02205748867102263, -0.04897121401367313, -0.07010672203707752, 0.29008673918970956, -0.0020876783821603432, 0.01606650988596511, -0.017521125982496902, -0.023198103633578725, -0.05776177307061119, -0.06120704146148565, -0.023766122331229312, -0.0005977689711279065, 0.0509316883364571, 0.010807225161789571, -0.027722531893364138, -0.05676363124281852, -0.030278860282686063, -0.003880991596474222]"
$ws.Range("B7").Value = 0.007442848126462526

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.1021468097139688
$ws.Range("B3").Value = 1.691963235174012
$ws.Range("B4").Value = 0.08775206191443111
$ws.Range("B5").Value = 0.9064504267509276
$ws.Range("B6").Value = "[1.0, -0.0011669597040826214, 0.015347392889375812, -0.011157823010246821, -0.047662370095994895, 0.20475804067996925, -0.0091473005465532, 0.000818845016198927, -0.01595360836876988, -0.0331062845154455, -0.053198552258009804, -0.05545262188146058, -0.02385264049525677, -0.00620782067620316, 0.038094902031899934, 0.009703227573229627, -0.029777991376613168, -0.04665223464333657, -0.01984237999413492, -0.008174620502161705]"
$ws.Range("B7").Value = -0.02055870320743971
$ws.Range("B9").Value = 0.802456396128738

$wb.Save()
